$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 564.36365
$ws.Range("I61").Value = 564.36365
$ws.Range("K61").Value = 1693.09095
$ws.Range("M61").Value = -1521.09095

$ws.Range("H116").Value = 4573.1055
$ws.Range("I116").Value = 1989.9
$ws.Range("J116").Value = 7443.3335
$ws.Range("K116").Value = 1989.9
$ws.Range("L116").Value = 7443.3335
$ws.Range("M116").Value = 1452.1
$ws.Range("N116").Value = -14327.3335

$ws.Range("H135").Value = 303608.28
$ws.Range("I135").Value = 323030.94
$ws.Range("K135").Value = 2907278.46
$ws.Range("M135").Value = -2904743.46

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2453032.8
$ws.Range("I32").Value = 2501954
$ws.Range("K32").Value = 2501954
$ws.Range("M32").Value = -2501667

$ws.Range("H110").Value = 1010.1429
$ws.Range("I110").Value = 934
$ws.Range("K110").Value = 934
$ws.Range("M110").Value = 1111

$ws.Range("H122").Value = 46999.8
$ws.Range("I122").Value = 103999.5
$ws.Range("J122").Value = 9000
$ws.Range("K122").Value = 311998.5
$ws.Range("L122").Value = 27000
$ws.Range("M122").Value = -309548.5
$ws.Range("N122").Value = -31900

$ws.Range("H132").Value = 2974.3137
$ws.Range("I132").Value = 1813.6177
$ws.Range("J132").Value = 5295.706
$ws.Range("K132").Value = 5440.8531
$ws.Range("L132").Value = 15887.118
$ws.Range("M132").Value = -2910.8531
$ws.Range("N132").Value = -20947.118

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2432.7144
$ws.Range("I94").Value = 1197.3334
$ws.Range("J94").Value = 4079.889
$ws.Range("K94").Value = 1197.3334
$ws.Range("L94").Value = 4079.889
$ws.Range("M94").Value = -746.3334
$ws.Range("N94").Value = -4981.889

$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws.Range("H105").Value = 41804.82
$ws.Range("I105").Value = 51779.773
$ws.Range("J105").Value = 3151.875
$ws.Range("K105").Value = 51779.773
$ws.Range("L105").Value = 3151.875
$ws.Range("M105").Value = -50032.773
$ws.Range("N105").Value = -6645.875

$ws.Range("H106").Value = 67825
$ws.Range("J106").Value = 67825
$ws.Range("L106").Value = 67825
$ws.Range("N106").Value = -70349

$ws.Range("H134").Value = 4162.551
$ws.Range("I134").Value = 3122.868
$ws.Range("K134").Value = 9368.603999999999
$ws.Range("M134").Value = -6833.603999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5241.5293
$ws.Range("J16").Value = 7714.2856
$ws.Range("L16").Value = 7714.2856
$ws.Range("N16").Value = -8288.285599999999

$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

$ws.Range("H99").Value = 9176.8125
$ws.Range("I99").Value = 10133.286
$ws.Range("K99").Value = 10133.286
$ws.Range("M99").Value = -8635.286

$ws.Range("H107").Value = 1879.7894
$ws.Range("I107").Value = 1576
$ws.Range("J107").Value = 2217.3333
$ws.Range("K107").Value = 1576
$ws.Range("L107").Value = 2217.3333
$ws.Range("M107").Value = 344
$ws.Range("N107").Value = -6057.3333

$ws.Range("H113").Value = 5241.5293
$ws.Range("J113").Value = 7714.2856
$ws.Range("L113").Value = 7714.2856
$ws.Range("N113").Value = -12054.2856

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws.Range("H126").Value = 9176.8125
$ws.Range("I126").Value = 10133.286
$ws.Range("K126").Value = 30399.858
$ws.Range("M126").Value = -27929.858

$ws.Range("H134").Value = 8893.091
$ws.Range("I134").Value = 9850.684999999999
$ws.Range("K134").Value = 29552.055
$ws.Range("M134").Value = -27017.055

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 111342.78
$ws.Range("J2").Value = 666691.7
$ws.Range("L2").Value = 4000150.2
$ws.Range("N2").Value = -4000376.2

$ws.Range("H34").Value = 3816.25
$ws.Range("J34").Value = 6131.5835
$ws.Range("L34").Value = 18394.7505
$ws.Range("N34").Value = -18562.7505

$ws.Range("H64").Value = 33360934
$ws.Range("J64").Value = 50040000
$ws.Range("L64").Value = 150120000
$ws.Range("N64").Value = -150120540

$ws.Range("H67").Value = 33360934
$ws.Range("J67").Value = 50040000
$ws.Range("L67").Value = 150120000
$ws.Range("N67").Value = -150121872

$ws.Range("H98").Value = 711.8
$ws.Range("I98").Value = 831.6667
$ws.Range("J98").Value = 532
$ws.Range("K98").Value = 2495.0001
$ws.Range("L98").Value = 1596
$ws.Range("M98").Value = -997.0001000000002
$ws.Range("N98").Value = -4592

$ws.Range("H107").Value = 33333910
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 33333910
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 100001730
$ws.Range("N107").Value = -100005570
$ws.Range("M107").ClearContents()

$ws.Range("H122").Value = 2831626.5
$ws.Range("J122").Value = 6249
$ws.Range("L122").Value = 56241
$ws.Range("N122").Value = -61141

$ws.Range("H131").Value = 1632.68
$ws.Range("J131").Value = 2745.2
$ws.Range("L131").Value = 8235.599999999999
$ws.Range("N131").Value = -18315.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 30000
$ws.Range("J48").Value = 30000
$ws.Range("L48").Value = 30000
$ws.Range("N48").Value = -30970

$ws.Range("H70").Value = 6424.609
$ws.Range("I70").Value = 5494.448
$ws.Range("K70").Value = 5494.448
$ws.Range("M70").Value = -5224.448

$ws.Range("H73").Value = 6424.609
$ws.Range("I73").Value = 5494.448
$ws.Range("K73").Value = 5494.448
$ws.Range("M73").Value = -4558.448

$ws.Range("H80").Value = 2170.2778
$ws.Range("I80").Value = 1848.125
$ws.Range("K80").Value = 1848.125
$ws.Range("M80").Value = -850.125

$ws.Range("H83").Value = 2170.2778
$ws.Range("I83").Value = 1848.125
$ws.Range("K83").Value = 9240.625
$ws.Range("M83").Value = -4248.625

$ws.Range("H102").Value = 9088.444
$ws.Range("I102").Value = 7633.1665
$ws.Range("J102").Value = 11999
$ws.Range("K102").Value = 7633.1665
$ws.Range("L102").Value = 11999
$ws.Range("M102").Value = -6011.1665
$ws.Range("N102").Value = -15243

$ws.Range("H122").Value = 43855.117
$ws.Range("I122").Value = 66853.06
$ws.Range("K122").Value = 200559.18
$ws.Range("M122").Value = -198109.18

$ws.Range("H132").Value = 3051.95
$ws.Range("I132").Value = 2986.4
$ws.Range("K132").Value = 8959.200000000001
$ws.Range("M132").Value = -6429.200000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 13237573
$ws.Range("J46").Value = 12347738
$ws.Range("L46").Value = 12347738
$ws.Range("N46").Value = -12348114

$ws.Range("H55").Value = 250000580
$ws.Range("I55").Value = 1000000000
$ws.Range("J55").Value = 771.3333
$ws.Range("K55").Value = 1000000000
$ws.Range("L55").Value = 771.3333
$ws.Range("M55").Value = -999999827
$ws.Range("N55").Value = -1117.3333

$ws.Range("H82").Value = 2364
$ws.Range("J82").Value = 10000
$ws.Range("L82").Value = 10000
$ws.Range("N82").Value = -10722

$ws.Range("H85").Value = 2364
$ws.Range("J85").Value = 10000
$ws.Range("L85").Value = 10000
$ws.Range("N85").Value = -12496

$ws.Range("H122").Value = 5739.885
$ws.Range("I122").Value = 5048.8
$ws.Range("K122").Value = 15146.4
$ws.Range("M122").Value = -12696.4

$ws.Range("H132").Value = 12826868
$ws.Range("I132").Value = 25005342
$ws.Range("J132").Value = 7421.263
$ws.Range("K132").Value = 75016026
$ws.Range("L132").Value = 22263.789
$ws.Range("M132").Value = -75013496
$ws.Range("N132").Value = -27323.789

$ws.Range("H136").Value = 7899.647
$ws.Range("I136").Value = 6073.5
$ws.Range("K136").Value = 18220.5
$ws.Range("M136").Value = -15670.5

$ws.Range("H139").Value = 88838.44500000001
$ws.Range("J139").Value = 88838.44500000001
$ws.Range("L139").Value = 88838.44500000001
$ws.Range("N139").Value = -99118.44500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H122").Value = 135616.5
$ws.Range("I122").Value = 201762.5
$ws.Range("K122").Value = 605287.5
$ws.Range("M122").Value = -602837.5

$ws.Range("H126").Value = 2453.3333
$ws.Range("I126").Value = 1528.8889
$ws.Range("K126").Value = 4586.6667
$ws.Range("M126").Value = -2116.6667

$ws.Range("H132").Value = 13524422
$ws.Range("I132").Value = 20003630
$ws.Range("J132").Value = 26073.834
$ws.Range("K132").Value = 60010890
$ws.Range("L132").Value = 78221.50199999999
$ws.Range("M132").Value = -60008360
$ws.Range("N132").Value = -83281.50199999999

$ws.Range("H136").Value = 62569560
$ws.Range("I136").Value = 250001230
$ws.Range("J136").Value = 92333.336
$ws.Range("K136").Value = 750003690
$ws.Range("L136").Value = 277000.008
$ws.Range("M136").Value = -750001140
$ws.Range("N136").Value = -282100.008
